# Generate Report for Handoff
# - Update status text "Handed back: in sync with en-US" -> "Ready for handoff"
#   on the Overview sheet (columns zh-cn/de-de) and on each language sheet's
#   "Status" column.
# - Update the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps
#   to reflect the new handoff generation time.
# - Shrink the width of the corresponding "Status" / zh-cn / de-de columns
#   to fit the shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# Target stored column width (per the OOXML <col> width attribute) is
# 17.2159881591797 characters. The COM ColumnWidth setter here quantizes to
# whole pixels (MDW=6 -> stored = (round(ColumnWidth*6)+5)/6), so
# 16.333333333333332 is the closest input that reproduces the intended
# (shrunk) column width.
$newWidth = 16.333333333333332

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-30 17:10:21"
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-30 17:10:03"
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-30 17:10:21"
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
